$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 15 and 16, pushing existing rows 15-23 down to 17-25
$ws.Rows("15:16").Insert()

# Row 15: emu / EMU membership / - / explanation
$ws.Range("A15").Value = "emu"
$ws.Range("B15").Value = "EMU membership"

# Row 16: reer / real effective exchange rate / (no sign) / explanation
$ws.Range("A16").Value = "reer"
$ws.Range("B16").Value = "real effective exchange rate"
$ws.Range("D16").Value = "not sure, yet to think about it"

# Fill in remaining row 15 cells last (matches shared-string ordering)
$ws.Range("D15").Value = "probably convergence in spreads should take place because investors wanna invest in euros"
$ws.Range("C15").Value = "-"

# Update the active selection to D15 as in the target file
$ws.Range("D15").Select()
